# Simulações para 55 Cnc e
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# First block (row 2) - update planet/star parameters
$ws.Range("C2").Value = 0.95
$ws.Range("D2").Value = 0.015440000000000001

# Second block (row 5) - update latitude + planet/star parameters
$ws.Range("B5").Value = 45
$ws.Range("C5").Value = 0.95
$ws.Range("D5").Value = 0.015440000000000001

# K10 mirrors the recalculated K5 result (static copy, not a formula)
$ws.Range("K10").Value = $ws.Range("K5").Value2
